$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '27.398.00'
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -0.71%  '
$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.639.73'
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -1.68%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '211.58'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -1.62%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.530'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +3.85%  '
$ws.Cells.Item(7, 5).Value = '  -0.04%  '
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '23.04'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -2.76%  '
$ws.Cells.Item(9, 5).Value = '  -2.10%  '
$ws.Cells.Item(10, 5).Value = '  -2.03%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0891'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +1.16%  '
$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.871.38'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -1.71%  '
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.628.99'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -2.56%  '
$ws.Cells.Item(14, 5).Value = '  -2.68%  '
$ws.Cells.Item(15, 5).Value = '  -0.43%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '64.34'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -2.84%  '
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '27.369.46'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.82%  '
$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '229.94'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -5.44%  '
$ws.Cells.Item(19, 5).Value = '  -1.20%  '
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.59'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -0.35%  '
$ws.Cells.Item(21, 5).Value = '  -0.01%  '
$ws.Cells.Item(22, 5).Value = '  -3.64%  '
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.52'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +2.44%  '
$ws.Cells.Item(24, 5).Value = '  -0.50%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '147.33'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.33%  '
$ws.Cells.Item(26, 5).Value = '  -3.16%  '
$ws.Cells.Item(27, 5).Value = '  +1.73%  '
$ws.Cells.Item(28, 5).Value = '  +0.13%  '
$ws.Cells.Item(29, 5).Value = '  -5.43%  '
$ws.Cells.Item(30, 5).Value = '  -3.78%  '
$ws.Cells.Item(31, 5).Value = '  -3.35%  '
$ws.Cells.Item(32, 5).Value = '  -2.09%  '
$ws.Cells.Item(33, 5).Value = '  +0.08%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.410.81'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -3.75%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.57'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +0.39%  '
$ws.Cells.Item(36, 5).Value = '  -0.31%  '
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.565'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -1.70%  '
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.881'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -5.00%  '
$ws.Cells.Item(39, 5).Value = '  -3.56%  '
$ws.Cells.Item(41, 5).Value = '  -0.03%  '
$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.52'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +1.86%  '
$ws.Cells.Item(43, 5).Value = '  -2.10%  '
$ws.Cells.Item(44, 5).Value = '  +0.53%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.792'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.48%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '64.49'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -7.20%  '
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.781.13'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -1.62%  '
$ws.Cells.Item(48, 5).Value = '  -4.54%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '87.43'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -2.15%  '
$ws.Cells.Item(50, 5).Value = '  -2.35%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0988'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -3.86%'
